$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "text" -> "Entry" for D2/D3
$ws.Range("D2").Value = "Entry"
$ws.Range("D3").Value = "Entry"

# "button" -> "Button" for D4/D5
$ws.Range("D4").Value = "Button"
$ws.Range("D5").Value = "Button"

# New row 6: log entry
$ws.Range("A6").Value = "log"
$ws.Range("B6").Value = "The Log"
$ws.Range("C6").Value = "Log box"
$ws.Range("D6").Value = "Text"

# Selection moved from B6 to A6
$ws.Range("A6").Select()
